$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Status" column header in F1, bold.
$ws.Range("F1").Value = "Status"
$ws.Range("F1").Font.Bold = $true

# Status values for the first several test cases.
$ws.Range("F2").Value = "unit test done"
$ws.Range("F3").Value = "unit test done"
$ws.Range("F4").Value = "unit test partially done"
$ws.Range("F5").Value = "unit test done"
$ws.Range("F6").Value = "unit test done"

# Move selection / view as in the saved file.
$ws.Range("F7").Select() | Out-Null
